$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.518.08"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "3.710.25"
$ws.Range("E3").Value = "  +8.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "3.699.37"
$ws.Range("E7").Value = "  +8.25%  "
$ws.Range("E8").Value = "  +4.40%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.201"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").Value = "4.306.55"
$ws.Range("E14").Value = "  +8.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "679.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").Value = "3.699.35"
$ws.Range("E17").Value = "  +7.89%  "
$ws.Range("D18").Value = "71.631.73"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +20.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.944"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.36%  "
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("E27").Value = "  +6.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "594.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.68%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "3.687.62"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").Value = "0.0₃0777"
$ws.Range("E40").Value = "  +7.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("E44").Value = "  +10.00%  "
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.83%  "
